# Fix "Recorded By" (column G) entries so "System" is no longer listed
# first/duplicated ahead of the real recorder's email/name.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "system, System, backup@backdoor.com"
#
# Applies to every row in the used range, on whichever sheet is active,
# matching by current cell text so the edit is robust to row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $t = $cell.Text
    if ($t -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    } elseif ($t -eq "System, system, backup@backdoor.com") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
